$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some Price values look like plain numbers to Excel (e.g. "696.90", "7.23")
# and would otherwise be auto-converted to numeric values, losing the exact
# text representation (trailing zeros, etc). Force those specific cells to
# text format first so the original string content is preserved exactly.
# (NumberFormat must be set per-cell; multi-area ranges only apply to the
# first area in this runtime.)
foreach ($addr in @("D5","D6","D11","D13","D18","D19","D21","D22","D23","D25","D27","D29","D31","D32","D33","D34","D35","D39","D40","D41","D45","D46","D47","D49","D50","D51")) {
    $ws.Range($addr).NumberFormat = "@"
}

# Column D = Price, Column E = Volume(1h)

$ws.Range("D2").Value = "71.092.64"
$ws.Range("E2").Value = "  +0.39%  "

$ws.Range("D3").Value = "3.865.65"
$ws.Range("E3").Value = "  +1.48%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "696.90"
$ws.Range("E5").Value = "  -0.57%  "

$ws.Range("D6").Value = "173.99"
$ws.Range("E6").Value = "  +0.37%  "

$ws.Range("D7").Value = "3.862.55"
$ws.Range("E7").Value = "  +1.45%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("E9").Value = "  +0.03%  "

$ws.Range("E10").Value = "  -0.07%  "

$ws.Range("D11").Value = "7.23"
$ws.Range("E11").Value = "  -5.83%  "

$ws.Range("E12").Value = "  -0.57%  "

$ws.Range("D13").Value = "0.0000259"
$ws.Range("E13").Value = "  +1.58%  "

$ws.Range("E14").Value = "  +0.88%  "

$ws.Range("D15").Value = "4.519.52"
$ws.Range("E15").Value = "  +1.56%  "

$ws.Range("D16").Value = "3.867.81"
$ws.Range("E16").Value = "  +1.55%  "

$ws.Range("D17").Value = "71.155.80"
$ws.Range("E17").Value = "  +0.42%  "

$ws.Range("D18").Value = "17.65"
$ws.Range("E18").Value = "  -0.34%  "

$ws.Range("D19").Value = "7.22"
$ws.Range("E19").Value = "  +0.18%  "

$ws.Range("E20").Value = "  -0.39%  "

$ws.Range("D21").Value = "11.01"
$ws.Range("E21").Value = "  -2.72%  "

$ws.Range("D22").Value = "497.69"
$ws.Range("E22").Value = "  +3.77%  "

$ws.Range("D23").Value = "0.722"
$ws.Range("E23").Value = "  +0.86%  "

$ws.Range("E24").Value = "  +3.18%  "

$ws.Range("D25").Value = "84.93"
$ws.Range("E25").Value = "  +1.33%  "

$ws.Range("E26").Value = "  +3.02%  "

$ws.Range("D27").Value = "12.23"
$ws.Range("E27").Value = "  -1.19%  "

$ws.Range("E28").Value = "  -0.64%  "

$ws.Range("D29").Value = "3.18"
$ws.Range("E29").Value = "  +2.43%  "

$ws.Range("E30").Value = "  +0.05%  "

$ws.Range("D31").Value = "7.58"
$ws.Range("E31").Value = "  +0.99%  "

$ws.Range("D32").Value = "2.29"
$ws.Range("E32").Value = "  -1.01%  "

$ws.Range("D33").Value = "29.69"
$ws.Range("E33").Value = "  +0.46%  "

$ws.Range("D34").Value = "0.183"
$ws.Range("E34").Value = "  +2.47%  "

$ws.Range("D35").Value = "9.23"
$ws.Range("E35").Value = "  +0.28%  "

$ws.Range("D36").Value = "3.818.78"
$ws.Range("E36").Value = "  +1.59%  "

$ws.Range("E37").Value = "  +0.15%  "

$ws.Range("E38").Value = "  +2.21%  "

$ws.Range("D39").Value = "2.39"
$ws.Range("E39").Value = "  +8.72%  "

# Row 40 and 41 swap (Mantle / dogwifhat)
$ws.Range("B40").Value = "Mantle"
$ws.Range("C40").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D40").Value = "1.05"
$ws.Range("E40").Value = "  +8.93%  "

$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").Value = "3.42"
$ws.Range("E41").Value = "  -3.71%  "

$ws.Range("E43").Value = "  +0.02%  "

$ws.Range("E44").Value = "  -0.03%  "

# Row 45 and 46 swap (FLOKI / Monero)
$ws.Range("B45").Value = "FLOKI"
$ws.Range("C45").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D45").Value = "0.000312"
$ws.Range("E45").Value = "  -6.89%  "

$ws.Range("B46").Value = "Monero"
$ws.Range("C46").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D46").Value = "163.68"
$ws.Range("E46").Value = "  +2.32%  "

$ws.Range("D47").Value = "49.43"
$ws.Range("E47").Value = "  +1.00%  "

$ws.Range("E48").Value = "  +1.13%  "

$ws.Range("D49").Value = "416.34"
$ws.Range("E49").Value = "  +3.50%  "

# Row 50 and 51 swap (ONDO / Arweave)
$ws.Range("B50").Value = "ONDO"
$ws.Range("C50").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D50").Value = "1.39"
$ws.Range("E50").Value = "  -1.90%  "

$ws.Range("B51").Value = "Arweave"
$ws.Range("C51").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D51").Value = "43.63"
$ws.Range("E51").Value = "  -4.23%  "
